# Fix model sheet -- 'type' now looks to match the prompt type and expands
# recursively. This removes the separate 'elementType' column on the model
# sheet and folds its value into the 'type' column directly, then makes the
# model sheet the active sheet/tab.

$wb = $excel.ActiveWorkbook

$model = $wb.Worksheets.Item("model")

# --- Update model sheet content ---
# Row 4: refrigerator_location  object -> geopoint (elementType column removed)
$model.Range("B4").Value = "geopoint"
# Row 5: refrigerator_condition  string -> select_one (matches the prompt/survey type)
$model.Range("B5").Value = "select_one"

# Remove the now-unused "elementType" column (column C) from the model sheet
$model.Columns.Item(3).Delete()

# --- Update selection / active sheet: "model" is now the active tab ---
$model.Range("C22").Select()
$model.Activate()

$wb.Save()
